$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.988.10"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  -2.28%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'1.797.93"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -0.04%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  +0.14%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'316.74"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  +1.09%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +0.13%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.5433"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'  +1.07%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.3799"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +0.65%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'0.07455"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  -0.98%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'41.90"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  -1.50%  "
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'1.092"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'  -1.98%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.16%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'6.214"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +0.74%  "
$ws.Range("E13").ClearFormats()
$ws.Range("B14").Value = "'Chainlink"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("C14").ClearFormats()
$ws.Range("D14").Value = "'7.399"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -0.18%  "
$ws.Range("E14").ClearFormats()
$ws.Range("B15").Value = "'Solana"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = "'https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("C15").ClearFormats()
$ws.Range("D15").Value = "'20.32"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'  -2.81%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'1.796.36"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  -0.06%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'89.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -1.59%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'0.00001063"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  -0.07%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.06516"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +1.19%  "
$ws.Range("E19").ClearFormats()
$ws.Range("B20").Value = "'Dai"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C20").ClearFormats()
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  -0.02%  "
$ws.Range("E20").ClearFormats()
$ws.Range("B21").Value = "'Avalanche"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").Value = "'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("C21").ClearFormats()
$ws.Range("D21").Value = "'17.40"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +0.93%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'5.935"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.11%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'28.030.10"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -2.16%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'11.15"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "'  -0.19%  "
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'2.092"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -0.73%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'156.38"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -2.64%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'20.34"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -0.55%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'2.004.28"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  +0.02%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.341"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  -1.64%  "
$ws.Range("E29").ClearFormats()
$ws.Range("D30").Value = "'122.38"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  -0.80%  "
$ws.Range("E30").ClearFormats()
$ws.Range("D31").Value = "'0.1105"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "'  +7.05%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +0.44%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'3.673"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -0.58%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'5.537"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  -2.13%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.06948"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +7.07%  "
$ws.Range("E35").ClearFormats()
$ws.Range("D36").Value = "'0.2205"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "'  -2.34%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.02289"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -1.15%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'5.080"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  +0.91%  "
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'8.442"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'  -4.92%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'11.24"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -0.86%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.6143"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -1.67%  "
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'1.169"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  -3.49%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  +1.62%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'13.31"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'  -0.42%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'3.685"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").Value = "'0.5728"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  -2.39%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'124.37"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -1.66%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'1.182"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +2.04%  "
$ws.Range("E48").ClearFormats()
$ws.Range("E49").Value = "'  -2.08%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.06808"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  -1.22%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.00000000295"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +38.65%  "
$ws.Range("E51").ClearFormats()
